$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phd")
$ws.Rows.Item(6).Select() | Out-Null
$ws.Rows.Item(6).Delete()
